$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 747; this shifts the existing rows 747:811 down to 748:812
# and pushes the dimension/used range to A1:R812.
$ws.Rows.Item(747).Insert()

# Populate the newly inserted row 747 with the new record's data.
$ws.Cells.Item(747, 1).Value = 3
$ws.Cells.Item(747, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(747, 3).Value = "Coquimbo"
$ws.Cells.Item(747, 4).Value = 45166
$ws.Cells.Item(747, 5).Value = 5
$ws.Cells.Item(747, 6).Value = 100112037
$ws.Cells.Item(747, 7).Value = "Cebollín"
$ws.Cells.Item(747, 8).Value = "Sin especificar"
$ws.Cells.Item(747, 9).Value = "Primera"
$ws.Cells.Item(747, 10).Value = 170
$ws.Cells.Item(747, 11).Value = 4000
$ws.Cells.Item(747, 12).Value = 4500
$ws.Cells.Item(747, 13).Value = 4176
$ws.Cells.Item(747, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(747, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(747, 16).Value = 116
$ws.Cells.Item(747, 17).Value = 36
$ws.Cells.Item(747, 18).Value = "Hortaliza"
